# Update gh-pages to output generated at 456a3b4
# Applies the "想去人数" (F column) count refresh across all four sheets,
# plus the 本地生活 sheet's G7 ticket-price cell flipping from a numeric
# price to the "不可售" (not for sale) status text.

$wb = $excel.ActiveWorkbook

# --- 展览 ("Exhibitions") ---------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 7746
$ws.Range("F6").Value  = 4862
$ws.Range("F7").Value  = 7152
$ws.Range("F9").Value  = 1522
$ws.Range("F10").Value = 878
$ws.Range("F11").Value = 211
$ws.Range("F12").Value = 67
$ws.Range("F13").Value = 1181
$ws.Range("F15").Value = 7
$ws.Range("F16").Value = 15
$ws.Range("F20").Value = 1219
$ws.Range("F24").Value = 1253
$ws.Range("F28").Value = 10
$ws.Range("F29").Value = 49
$ws.Range("F30").Value = 213
$ws.Range("F33").Value = 10
$ws.Range("F34").Value = 134
$ws.Range("F35").Value = 131
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 565
$ws.Range("F40").Value = 87
$ws.Range("F41").Value = 63
$ws.Range("F42").Value = 102
$ws.Range("F43").Value = 418
$ws.Range("F44").Value = 1205
$ws.Range("F45").Value = 599
$ws.Range("F46").Value = 158
$ws.Range("F47").Value = 28
$ws.Range("F48").Value = 26

# --- 演出 ("Shows") -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 32
$ws.Range("F13").Value = 12
$ws.Range("F15").Value = 1735
$ws.Range("F29").Value = 10
$ws.Range("F32").Value = 872
$ws.Range("F34").Value = 996
$ws.Range("F35").Value = 616
$ws.Range("F42").Value = 147

# --- 本地生活 ("Local life") --------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 736
$ws.Range("F5").Value  = 861
$ws.Range("F6").Value  = 695
# F7 (293) is unchanged; G7 flips from a numeric price (10) to status text.
$ws.Range("G7").Value  = "不可售"
$ws.Range("F8").Value  = 152
$ws.Range("F9").Value  = 97
$ws.Range("F10").Value = 1710
$ws.Range("F11").Value = 2615

# --- 全部类型 ("All types") ---------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 736
$ws.Range("F4").Value  = 861
$ws.Range("F7").Value  = 695
$ws.Range("F8").Value  = 695
$ws.Range("F9").Value  = 7746
$ws.Range("F10").Value = 152
$ws.Range("F11").Value = 4862
$ws.Range("F12").Value = 7152
$ws.Range("F14").Value = 1522
$ws.Range("F15").Value = 878
$ws.Range("F16").Value = 97
$ws.Range("F17").Value = 211
$ws.Range("F18").Value = 1710
$ws.Range("F19").Value = 2615
$ws.Range("F21").Value = 67
$ws.Range("F22").Value = 1181
$ws.Range("F25").Value = 1219
$ws.Range("F28").Value = 1253
$ws.Range("F30").Value = 213
$ws.Range("F32").Value = 10
$ws.Range("F33").Value = 872
$ws.Range("F34").Value = 134
$ws.Range("F35").Value = 131
$ws.Range("F36").Value = 996
$ws.Range("F37").Value = 565
$ws.Range("F38").Value = 616
$ws.Range("F39").Value = 87
$ws.Range("F40").Value = 102
$ws.Range("F42").Value = 418
$ws.Range("F43").Value = 599
$ws.Range("F45").Value = 147
$ws.Range("F46").Value = 158
$ws.Range("F49").Value = 28
